$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '65.962.40'
$ws.Range('E2').Value = '  -2.73%  '

# Row 3
$ws.Range('D3').Value = '3.478.19'
$ws.Range('E3').Value = '  +0.58%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.33'
$ws.Range('E5').Value = '  -1.52%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.67'
$ws.Range('E6').Value = '  -3.71%  '

# Row 7
$ws.Range('E7').Value = '  +0.01%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.596'
$ws.Range('E8').Value = '  -3.02%  '

# Row 9
$ws.Range('D9').Value = '3.479.25'
$ws.Range('E9').Value = '  +0.73%  '

# Row 10
$ws.Range('E10').Value = '  -5.63%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.85'
$ws.Range('E11').Value = '  -1.90%  '

# Row 12
$ws.Range('E12').Value = '  -4.10%  '

# Row 13
$ws.Range('D13').Value = '4.082.60'
$ws.Range('E13').Value = '  +0.69%  '

# Row 14
$ws.Range('E14').Value = '  +0.19%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '30.06'
$ws.Range('E15').Value = '  -6.12%  '

# Row 16
$ws.Range('D16').Value = '66.058.53'
$ws.Range('E16').Value = '  -2.52%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000171'
$ws.Range('E17').Value = '  -3.08%  '

# Row 18
$ws.Range('D18').Value = '3.480.94'
$ws.Range('E18').Value = '  +0.67%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.92'
$ws.Range('E19').Value = '  -4.00%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.91'
$ws.Range('E20').Value = '  -1.16%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '366.43'
$ws.Range('E21').Value = '  -6.29%  '

# Row 22
$ws.Range('E22').Value = '  -1.68%  '

# Row 23
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.57'
$ws.Range('E23').Value = '  +1.16%  '

# Row 24
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.03%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.535'
$ws.Range('E25').Value = '  -0.36%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000124'
$ws.Range('E26').Value = '  +3.63%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.60'
$ws.Range('E27').Value = '  -7.54%  '

# Row 28
$ws.Range('E28').Value = '  +1.24%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.06%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '24.02'
$ws.Range('E30').Value = '  +2.50%  '

# Row 31
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.99'
$ws.Range('E31').Value = '  -2.92%  '

# Row 32
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.76'
$ws.Range('E32').Value = '  -5.55%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.04%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.14'
$ws.Range('E34').Value = '  -2.63%  '

# Row 35
$ws.Range('E35').Value = '  -7.62%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.53'
$ws.Range('E36').Value = '  -1.65%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.98'
$ws.Range('E37').Value = '  -2.18%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '29.13'
$ws.Range('E38').Value = '  +11.88%  '

# Row 39
$ws.Range('E39').Value = '  +0.33%  '

# Row 40
$ws.Range('D40').Value = '2.828.20'
$ws.Range('E40').Value = '  +4.00%  '

# Row 41
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.76'
$ws.Range('E41').Value = '  -5.58%  '

# Row 42
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.59'
$ws.Range('E42').Value = '  -6.67%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.45'
$ws.Range('E43').Value = '  -3.32%  '

# Row 44
$ws.Range('E44').Value = '  -3.89%  '

# Row 45
$ws.Range('E45').Value = '  -4.87%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.92'
$ws.Range('E46').Value = '  -3.28%  '

# Row 47
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0289'
$ws.Range('E47').Value = '  -2.84%  '

# Row 48
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '24.07'
$ws.Range('E48').Value = '  -7.99%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '309.60'
$ws.Range('E49').Value = '  -5.70%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.820'
$ws.Range('E50').Value = '  -2.34%  '

# Row 51
$ws.Range('B51').Value = 'ONDO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.980'
$ws.Range('E51').Value = '  -5.91%  '
